# Refresh the Price (column D) and Volume(1h) (column E) figures for the crypto
# ranking snapshot -- values taken from the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Most new values (percentages, and prices that contain more than one "." like
# "52.245.21") are never mistaken for a number by Excel, so they can be written
# directly as plain text.
$plainUpdates = @{
    'D2' = '52.245.21'
    'E2' = '  -0.03%  '
    'D3' = '2.824.08'
    'E3' = '  +0.95%  '
    'E4' = '  +0.02%  '
    'E5' = '  +2.64%  '
    'E6' = '  -3.62%  '
    'E7' = '  +3.38%  '
    'E8' = '  +0.08%  '
    'E9' = '  +1.18%  '
    'E10' = '  -4.81%  '
    'E11' = '  +0.89%  '
    'E12' = '  +0.73%  '
    'E13' = '  -0.90%  '
    'D15' = '3.266.62'
    'E15' = '  +0.97%  '
    'D16' = '2.823.86'
    'E16' = '  +0.92%  '
    'E17' = '  +3.91%  '
    'D18' = '52.092.88'
    'E18' = '  -0.15%  '
    'E19' = '  +4.84%  '
    'E20' = '  -1.29%  '
    'E21' = '  -0.09%  '
    'E22' = '  +1.75%  '
    'E23' = '  +0.85%  '
    'E24' = '  +0.55%  '
    'E25' = '  +2.36%  '
    'E26' = '  +1.02%  '
    'E27' = '  +0.14%  '
    'E28' = '  +1.12%  '
    'E29' = '  +0.51%  '
    'E31' = '  +10.86%  '
    'E32' = '  +5.01%  '
    'E33' = '  +0.03%  '
    'E34' = '  +4.19%  '
    'E35' = '  +12.08%  '
    'E36' = '  +3.47%  '
    'E37' = '  +0.01%  '
    'E38' = '  +1.79%  '
    'E39' = '  -3.70%  '
    'E40' = '  -2.79%  '
    'E41' = '  +1.63%  '
    'E42' = '  +0.09%  '
    'E43' = '  -5.87%  '
    'E44' = '  -0.76%  '
    'E45' = '  -1.93%  '
    'E46' = '  +0.41%  '
    'D47' = '2.087.88'
    'E47' = '  +0.76%  '
    'E48' = '  -4.47%  '
    'E49' = '  +7.58%  '
    'E50' = '  +1.09%  '
    'E51' = '  +2.63%  '
}
foreach ($cell in $plainUpdates.Keys) {
    $ws.Range($cell).Value = $plainUpdates[$cell]
}

# Some new prices look like an ordinary decimal (e.g. "356.22", "41.00", "0.0482").
# Assigning those to .Value directly would make Excel coerce them to a Double,
# silently dropping trailing zeros / introducing floating-point noise ("41" instead
# of "41.00"). A leading apostrophe forces Excel to keep them as text, exactly like
# the source data; re-applying the Normal style afterwards keeps the cell format
# identical to the rest of the column.
$textUpdates = @{
    'D5' = "'356.22"
    'D6' = "'112.15"
    'D7' = "'0.570"
    'D10' = "'41.00"
    'D13' = "'19.96"
    'D14' = "'7.77"
    'D17' = "'0.930"
    'D19' = "'7.50"
    'D20' = "'3.21"
    'D21' = "'13.44"
    'D23' = "'70.75"
    'D24' = "'271.64"
    'D25' = "'2.84"
    'D26' = "'27.00"
    'D28' = "'10.34"
    'D31' = "'0.0482"
    'D32' = "'52.77"
    'D34' = "'5.95"
    'D36' = "'0.0855"
    'D38' = "'3.28"
    'D40' = "'18.38"
    'D42' = "'127.69"
    'D44' = "'23.23"
    'D46' = "'3.36"
    'D50' = "'0.982"
    'D51' = "'9.20"
}
foreach ($cell in $textUpdates.Keys) {
    $rng = $ws.Range($cell)
    $rng.Value = $textUpdates[$cell]
    $rng.Style = "Normal"
}
